# Auto-generated script applying Leve profit recalculation updates
# to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets of the workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 9300
$ws.Range("I18").Value = 3933.3333
$ws.Range("J18").Value = 11312.5
$ws.Range("K18").Value = 3933.3333
$ws.Range("L18").Value = 11312.5
$ws.Range("M18").Value = -3649.3333
$ws.Range("N18").Value = -11880.5
$ws.Range("H33").Value = 1078827.2
$ws.Range("I33").Value = 1232617.2
$ws.Range("J33").Value = 2297.5
$ws.Range("K33").Value = 1232617.2
$ws.Range("L33").Value = 2297.5
$ws.Range("M33").Value = -1232388.2
$ws.Range("N33").Value = -2755.5
$ws.Range("H40").Value = 1905
$ws.Range("I40").Value = 1144
$ws.Range("J40").Value = 2181.7273
$ws.Range("K40").Value = 1144
$ws.Range("L40").Value = 2181.7273
$ws.Range("M40").Value = -969
$ws.Range("N40").Value = -2531.7273
$ws.Range("H58").Value = 484.41666
$ws.Range("J58").Value = 1700
$ws.Range("L58").Value = 5100
$ws.Range("N58").Value = -5400
$ws.Range("H69").Value = 3500
$ws.Range("J69").Value = 3500
$ws.Range("L69").Value = 10500
$ws.Range("N69").Value = -12248
$ws.Range("H72").Value = 3500
$ws.Range("J72").Value = 3500
$ws.Range("L72").Value = 31500
$ws.Range("N72").Value = -40236
$ws.Range("H112").Value = 1875.8
$ws.Range("J112").Value = 1880.3684
$ws.Range("L112").Value = 5641.1052
$ws.Range("N112").Value = -7857.1052
$ws.Range("H137").Value = 51681.55
$ws.Range("I137").Value = 1782.5333
$ws.Range("K137").Value = 5347.5999
$ws.Range("M137").Value = -2797.5999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4497.8
$ws.Range("I2").Value = 3998
$ws.Range("K2").Value = 3998
$ws.Range("M2").Value = -3885
$ws.Range("H32").Value = 144680.2
$ws.Range("I32").Value = 137975.73
$ws.Range("K32").Value = 137975.73
$ws.Range("M32").Value = -137688.73
$ws.Range("H61").Value = 3301.5
$ws.Range("J61").Value = 2889.5
$ws.Range("L61").Value = 2889.5
$ws.Range("N61").Value = -3313.5
$ws.Range("H74").Value = 1184
$ws.Range("I74").Value = 903.7646999999999
$ws.Range("K74").Value = 903.7646999999999
$ws.Range("M74").Value = -29.76469999999995
$ws.Range("H77").Value = 1184
$ws.Range("I77").Value = 903.7646999999999
$ws.Range("K77").Value = 4518.8235
$ws.Range("M77").Value = -150.8234999999995
$ws.Range("H116").Value = 4497.8
$ws.Range("I116").Value = 3998
$ws.Range("K116").Value = 3998
$ws.Range("M116").Value = -1704
$ws.Range("H132").Value = 1571.0526
$ws.Range("I132").Value = 1578.1875
$ws.Range("J132").Value = 1533
$ws.Range("K132").Value = 4734.5625
$ws.Range("L132").Value = 4599
$ws.Range("M132").Value = -2204.5625
$ws.Range("N132").Value = -9659
$ws.Range("H136").Value = 3301.5
$ws.Range("J136").Value = 2889.5
$ws.Range("L136").Value = 8668.5
$ws.Range("N136").Value = -13768.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4497.8
$ws.Range("I3").Value = 3998
$ws.Range("K3").Value = 3998
$ws.Range("M3").Value = -3884
$ws.Range("H20").Value = 9393.182000000001
$ws.Range("I20").Value = 8238.941000000001
$ws.Range("K20").Value = 8238.941000000001
$ws.Range("M20").Value = -7991.941000000001
$ws.Range("H134").Value = 1812.2
$ws.Range("I134").Value = 1831.5
$ws.Range("K134").Value = 5494.5
$ws.Range("M134").Value = -2959.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 4499.5
$ws.Range("J8").Value = 3999
$ws.Range("L8").Value = 3999
$ws.Range("N8").Value = -4279
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H31").Value = 2146.4167
$ws.Range("I31").Value = 1950.6364
$ws.Range("K31").Value = 1950.6364
$ws.Range("M31").Value = -1655.6364
$ws.Range("H34").Value = 2146.4167
$ws.Range("I34").Value = 1950.6364
$ws.Range("K34").Value = 1950.6364
$ws.Range("M34").Value = -1748.6364
$ws.Range("H50").Value = 14998.333
$ws.Range("J50").Value = 14998.333
$ws.Range("L50").Value = 14998.333
$ws.Range("N50").Value = -16248.333
$ws.Range("H51").Value = 14998.143
$ws.Range("J51").Value = 14998.143
$ws.Range("L51").Value = 14998.143
$ws.Range("N51").Value = -16470.143
$ws.Range("H58").Value = 1306
$ws.Range("H60").Value = 13297.714
$ws.Range("J60").Value = 14998.2
$ws.Range("L60").Value = 14998.2
$ws.Range("N60").Value = -16020.2
$ws.Range("H61").Value = 14998.143
$ws.Range("J61").Value = 14998.143
$ws.Range("L61").Value = 14998.143
$ws.Range("N61").Value = -15694.143
$ws.Range("H68").Value = 22857.143
$ws.Range("H71").Value = 22857.143
$ws.Range("H107").Value = 910.0769
$ws.Range("I107").Value = 453.2857
$ws.Range("K107").Value = 453.2857
$ws.Range("H132").Value = 2738.6667
$ws.Range("I132").Value = 2524.75
$ws.Range("J132").Value = 3166.5
$ws.Range("K132").Value = 7574.25
$ws.Range("L132").Value = 9499.5
$ws.Range("M132").Value = -5044.25
$ws.Range("N132").Value = -14559.5
$ws.Range("H136").Value = 1306

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 3405.3333
$ws.Range("J118").Value = 4048.9
$ws.Range("L118").Value = 12146.7
$ws.Range("N118").Value = -14632.7
$ws.Range("H131").Value = 2210.625
$ws.Range("J131").Value = 2449.6365
$ws.Range("L131").Value = 7348.9095
$ws.Range("N131").Value = -17428.9095

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 30472120
$ws.Range("I11").Value = 5515322
$ws.Range("J11").Value = 80385710
$ws.Range("K11").Value = 5515322
$ws.Range("L11").Value = 80385710
$ws.Range("M11").Value = -5515183
$ws.Range("N11").Value = -80385988
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H18").Value = 30000
$ws.Range("I18").Value = 30000
$ws.Range("K18").Value = 30000
$ws.Range("M18").Value = -29707
$ws.Range("H20").Value = 34241.637
$ws.Range("I20").Value = 4366.6665
$ws.Range("J20").Value = 45444.75
$ws.Range("K20").Value = 4366.6665
$ws.Range("L20").Value = 45444.75
$ws.Range("M20").Value = -4121.6665
$ws.Range("N20").Value = -45934.75
$ws.Range("H126").Value = 3338.5715
$ws.Range("I126").Value = 2274
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 6822
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -4352
$ws.Range("N126").Value = -22940
$ws.Range("H131").Value = 63494.5
$ws.Range("J131").Value = 63494.5
$ws.Range("L131").Value = 63494.5
$ws.Range("N131").Value = -73574.5
$ws.Range("H132").Value = 2108.7778
$ws.Range("I132").Value = 2042.5834
$ws.Range("J132").Value = 2241.1667
$ws.Range("K132").Value = 6127.7502
$ws.Range("L132").Value = 6723.500100000001
$ws.Range("M132").Value = -3597.7502
$ws.Range("N132").Value = -11783.5001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1100.1818
$ws.Range("I22").Value = 1170
$ws.Range("J22").Value = 402
$ws.Range("K22").Value = 1170
$ws.Range("L22").Value = 402
$ws.Range("M22").Value = -875
$ws.Range("N22").Value = -992
$ws.Range("H27").Value = 1100.1818
$ws.Range("I27").Value = 1170
$ws.Range("J27").Value = 402
$ws.Range("K27").Value = 1170
$ws.Range("L27").Value = 402
$ws.Range("M27").Value = -1063
$ws.Range("N27").Value = -616
$ws.Range("H108").Value = 733542
$ws.Range("J108").Value = 733542
$ws.Range("L108").Value = 733542
$ws.Range("N108").Value = -741222
$ws.Range("H132").Value = 2561.875
$ws.Range("I132").Value = 1928.2142
$ws.Range("K132").Value = 5784.642599999999
$ws.Range("M132").Value = -3254.642599999999
$ws.Range("H136").Value = 4955.4287
$ws.Range("I136").Value = 2674.75
$ws.Range("K136").Value = 8024.25
$ws.Range("M136").Value = -5474.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 7500
$ws.Range("J28").Value = 7500
$ws.Range("L28").Value = 7500
$ws.Range("N28").Value = -8196
$ws.Range("H31").Value = 12500
$ws.Range("J31").Value = 12500
$ws.Range("L31").Value = 12500
$ws.Range("N31").Value = -13196
$ws.Range("H49").Value = 12250
$ws.Range("J49").Value = 19500
$ws.Range("L49").Value = 19500
$ws.Range("N49").Value = -19960
$ws.Range("H81").Value = 2392.1
$ws.Range("J81").Value = 1933
$ws.Range("L81").Value = 3866
$ws.Range("N81").Value = -5988
$ws.Range("H84").Value = 2392.1
$ws.Range("J84").Value = 1933
$ws.Range("L84").Value = 19330
$ws.Range("N84").Value = -29938
$ws.Range("H126").Value = 2639.652
$ws.Range("I126").Value = 2492.25
$ws.Range("J126").Value = 2976.5715
$ws.Range("K126").Value = 7476.75
$ws.Range("L126").Value = 8929.7145
$ws.Range("M126").Value = -5006.75
$ws.Range("N126").Value = -13869.7145
$ws.Range("H132").Value = 11452.588
$ws.Range("I132").Value = 16481.363
$ws.Range("K132").Value = 49444.08900000001
$ws.Range("M132").Value = -46914.08900000001

Write-Host "Applied all Leve profit updates"
